$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "60.887.68"
$c.Style = $s
$c = $ws.Range("E2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.53%  "
$c.Style = $s

# Row 3: Ethereum
$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.420.18"
$c.Style = $s
$c = $ws.Range("E3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.47%  "
$c.Style = $s

# Row 4: TetherUSD
$c = $ws.Range("D4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.994"
$c.Style = $s
$c = $ws.Range("E4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.47%  "
$c.Style = $s

# Row 5: BNB
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "570.86"
$c.Style = $s
$c = $ws.Range("E5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.58%  "
$c.Style = $s

# Row 6: Solana
$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "139.87"
$c.Style = $s
$c = $ws.Range("E6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.77%  "
$c.Style = $s

# Row 7: USDC
$c = $ws.Range("E7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.13%  "
$c.Style = $s

# Row 8: XRP
$c = $ws.Range("E8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.34%  "
$c.Style = $s

# Row 9: LidoStakedEther
$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.404.20"
$c.Style = $s
$c = $ws.Range("E9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.98%  "
$c.Style = $s

# Row 10: Dogecoin
$c = $ws.Range("E10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.24%  "
$c.Style = $s

# Row 11: TRON
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.160"
$c.Style = $s
$c = $ws.Range("E11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.24%  "
$c.Style = $s

# Row 12: Toncoin
$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.10"
$c.Style = $s
$c = $ws.Range("E12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.29%  "
$c.Style = $s

# Row 13: Cardano
$c = $ws.Range("E13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.45%  "
$c.Style = $s

# Row 14: Avalanche
$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "25.94"
$c.Style = $s
$c = $ws.Range("E14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.14%  "
$c.Style = $s

# Row 15: ShibaInu
$c = $ws.Range("E15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -4.06%  "
$c.Style = $s

# Row 16: WrappedliquidstakedEther2.0
$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.829.83"
$c.Style = $s
$c = $ws.Range("E16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.35%  "
$c.Style = $s

# Row 17: WrappedBTC
$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "60.726.98"
$c.Style = $s
$c = $ws.Range("E17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.31%  "
$c.Style = $s

# Row 18: WrappedEther
$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.395.05"
$c.Style = $s
$c = $ws.Range("E18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.44%  "
$c.Style = $s

# Row 19: Uniswap
$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.53"
$c.Style = $s
$c = $ws.Range("E19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +5.44%  "
$c.Style = $s

# Row 20: Chainlink
$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.68"
$c.Style = $s
$c = $ws.Range("E20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.99%  "
$c.Style = $s

# Row 21: BitcoinCash
$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "322.35"
$c.Style = $s
$c = $ws.Range("E21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.96%  "
$c.Style = $s

# Row 22: Polkadot
$c = $ws.Range("E22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.05%  "
$c.Style = $s

# Row 23: LEO
$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.06"
$c.Style = $s
$c = $ws.Range("E23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.00%  "
$c.Style = $s

# Row 24: Dai
$c = $ws.Range("E24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "
$c.Style = $s

# Row 25: SuiNetwork
$c = $ws.Range("E25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -5.43%  "
$c.Style = $s

# Row 26: Litecoin
$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "64.82"
$c.Style = $s
$c = $ws.Range("E26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.37%  "
$c.Style = $s

# Row 27: Bittensor
$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "582.96"
$c.Style = $s
$c = $ws.Range("E27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.33%  "
$c.Style = $s

# Row 28: Aptos
$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.43"
$c.Style = $s
$c = $ws.Range("E28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -9.96%  "
$c.Style = $s

# Row 29: WrappedeETH
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.514.93"
$c.Style = $s
$c = $ws.Range("E29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.38%  "
$c.Style = $s

# Row 30: PEPE
$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0₃0923"
$c.Style = $s
$c = $ws.Range("E30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -4.91%  "
$c.Style = $s

# Row 31: InternetComputer(DFINITY)
$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.91"
$c.Style = $s
$c = $ws.Range("E31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.90%  "
$c.Style = $s

# Row 32: Fetch.AI
$c = $ws.Range("E32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -5.97%  "
$c.Style = $s

# Row 33: PancakeSwap
$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = $s
$c = $ws.Range("E33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.96%  "
$c.Style = $s

# Row 34: Kaspa
$c = $ws.Range("E34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.85%  "
$c.Style = $s

# Row 35: FirstDigitalUSD
$c = $ws.Range("E35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.Style = $s

# Row 36: NEARProtocol
$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.63"
$c.Style = $s
$c = $ws.Range("E36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -6.44%  "
$c.Style = $s

# Row 37: ImmutableX
$c = $ws.Range("E37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -3.68%  "
$c.Style = $s

# Row 38: Monero
$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "151.02"
$c.Style = $s
$c = $ws.Range("E38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.57%  "
$c.Style = $s

# Row 39: PolygonEcosystemToken
$c = $ws.Range("E39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -3.14%  "
$c.Style = $s

# Row 40: EthereumClassic
$c = $ws.Range("E40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.02%  "
$c.Style = $s

# Row 41: RenderToken
$c = $ws.Range("E41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -4.04%  "
$c.Style = $s

# Row 42: USDe
$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $s
$c = $ws.Range("E42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.10%  "
$c.Style = $s

# Row 43: Stacks
$c = $ws.Range("E43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -3.27%  "
$c.Style = $s

# Row 44: OKB
$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "41.18"
$c.Style = $s
$c = $ws.Range("E44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -4.53%  "
$c.Style = $s

# Row 45: dogwifhat
$c = $ws.Range("E45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -4.88%  "
$c.Style = $s

# Row 46: BabyDogeCoin
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0₆0285"
$c.Style = $s
$c = $ws.Range("E46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +9.26%  "
$c.Style = $s

# Row 47: Aave
$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "141.05"
$c.Style = $s
$c = $ws.Range("E47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.61%  "
$c.Style = $s

# Row 48: Filecoin
$c = $ws.Range("E48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -4.28%  "
$c.Style = $s

# Row 49: Mantle
$c = $ws.Range("E49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -3.21%  "
$c.Style = $s

# Row 50: InjectiveProtocol
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "19.55"
$c.Style = $s
$c = $ws.Range("E50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.00%  "
$c.Style = $s

# Row 51: Hedera
$c = $ws.Range("E51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -4.26%  "
$c.Style = $s
